$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New A1:A33 values (column of returns data) replacing the prior sample.
$values = @(
    0.09611372762975634,
    -0.0059999999149979999,
    -0.0039999999354982663,
    -0.0079999998792477101,
    -0.0029999999491954199,
    0.082224276418463305,
    -0.0099999998562489978,
    -0.0099999998572721793,
    -0.0019999999740418772,
    -0.0019999999797164492,
    -0.0029999999661507459,
    -0.003499999960530964,
    0.027208633451657249,
    -0.0079999999091597829,
    -0.0010000000075889304,
    0.0091974396826368654,
    -0.0019999999939441793,
    -0.0039999999657345242,
    -0.006488739495079443,
    -0.0039999999377648976,
    -0.0039999999370881056,
    -0.0039999999365107897,
    -0.0049999999166931985,
    -0.019999999699384041,
    -0.019999999695041737,
    -0.0024999999563313224,
    -0.0024999999523855898,
    -0.0019999999380804212,
    -0.0069999998564753696,
    -0.059999999109968893,
    0.009282035360714147,
    -0.0099999998056947703,
    -0.0039999998892401578
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

# Column A narrowed by one character (16.42578125 -> 15.42578125 stored width).
$ws.Columns.Item(1).ColumnWidth = 14.592447916666666
